$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colC = @(0.3835736843423803, 0.381567749823148, 0.3805053764084505, 0.380115090030344, 0.3800528607501548, 0.3804999401258158, 0.3828469465027098, 0.388790032304243, 0.3939712749473472, 0.3965048554606483, 0.3974895968494536, 0.3972763896909726, 0.3965853634661869, 0.3961653866252277, 0.3938092475993642, 0.3924090185319073, 0.3916202708644505, 0.3913560722026403, 0.3925563549514948, 0.3967876477136372, 0.3997006509377741, 0.3981324419226553, 0.3924896934957616, 0.3870389998504322)
$colD = @(0.1168614100633079, 0.11480254448864, 0.1135653001970809, 0.1130679333871782, 0.1129857595392139, 0.1135585648467838, 0.1161459601371817, 0.1214312862328697, 0.1254410610454499, 0.1272922697124415, 0.1279971309215, 0.1278451561380081, 0.1273501822443421, 0.1270474960488173, 0.125320621716682, 0.1242681531479377, 0.1236653590140691, 0.1234617042710084, 0.1243799258816267, 0.1274954640543342, 0.1295540564549498, 0.1284533144210229, 0.1243293863110466, 0.1199790638178229)
$colE = @(0.1538977088172029, 0.1542561673580938, 0.1545286047948196, 0.1546528119365824, 0.1546742336100966, 0.1545302264704453, 0.15401045628162, 0.1534054988562232, 0.1532123858528287, 0.1531788629463833, 0.1531739582996146, 0.1531746685028814, 0.1531783034461185, 0.1531815437547408, 0.1532156679187686, 0.1532505001653703, 0.1532756494090464, 0.1532850436351119, 0.1532462630263858, 0.1531770245394135, 0.1531771676411573, 0.1531729449470518, 0.1532481626755633, 0.1535249283378128)
$colF = @(2.281331322514944, 2.29982489110462, 2.312444690789896, 2.317905275092322, 2.318831199265063, 2.31251704705663, 2.287445361374559, 2.248319706786887, 2.225702992611275, 2.216746811201432, 2.213547070970492, 2.214227660140466, 2.216479722164735, 2.217884155163404, 2.226315125292047, 2.231828616321721, 2.23512521908269, 2.236262921912569, 2.231228717152149, 2.215813030808761, 2.206855867874211, 2.21153412479795, 2.231499536526044, 2.25782853809541)
$colG = @(0.002537459553084928, 0.002540781978211733, 0.002542930364474509, 0.002543833190901006, 0.002543984758374019, 0.002542942429490896, 0.002538582676093545, 0.002530889558617799, 0.002525754159455532, 0.002523529004944576, 0.002522702267703664, 0.002522879615239004, 0.002523460670989199, 0.002523818650150737, 0.002525901805068639, 0.002527208119521465, 0.002527969925900488, 0.002528229657450146, 0.002527067979079857, 0.002523289570303806, 0.00252091269072167, 0.00252217283472429, 0.002527131302887176, 0.00253287962484069)
$colJ = @(0.227037872542077, 0.2290034245979946, 0.2302973117211433, 0.2308464803060737, 0.2309389924135186, 0.2303046293131636, 0.2276975406358765, 0.2232749232833697, 0.2204452521266624, 0.2192488826062906, 0.2188089010795409, 0.2189030783602348, 0.2192124234128343, 0.2194036064400677, 0.2205252653604219, 0.2212366328393074, 0.2216543454464563, 0.2217972449714587, 0.2211600214288332, 0.2191212069561743, 0.217864830107537, 0.2185284210157867, 0.2211946302282115, 0.2243976004365891)
$colK = @(1.724842512487612, 1.583822882290406, 1.497314805546239, 1.462083964795482, 1.456235297272372, 1.496839578551914, 1.676203880741411, 2.028484022625832, 2.287563225648512, 2.405467550044193, 2.450120031846836, 2.440503148088794, 2.409141055528437, 2.389931411585167, 2.279858674634909, 2.212343298766882, 2.173514879058757, 2.160369107941676, 2.219529961818125, 2.418352753029183, 2.548321170182817, 2.478952904519304, 2.216280911567026, 1.933132143756154)
$colL = @(0.1326414652375352, 0.1332869932048872, 0.1337116035587869, 0.1338917552065144, 0.133922099687819, 0.1337140043006411, 0.1328581882074964, 0.1314034488902625, 0.1304699988953679, 0.130074544809343, 0.1299289773561245, 0.1299601421064835, 0.1300624851405416, 0.1301257174896016, 0.1304964287360324, 0.1307313115936317, 0.130869157061035, 0.1309163013818466, 0.1307060236981101, 0.1300323111050954, 0.1296163733083251, 0.1298361414032154, 0.1307174476074895, 0.1317731604657357)
$colM = @(0.5589997316786253, 0.5321245139767399, 0.5157335853135692, 0.5090823984674842, 0.5079796934402978, 0.515643770103118, 0.5497104672487225, 0.6173765319759497, 0.6676000882677329, 0.6905557374432476, 0.6992637180822641, 0.6973876309723437, 0.6912718467290802, 0.6875277149930383, 0.6661020308070391, 0.6529855820577666, 0.6454516061258389, 0.6429025072276104, 0.6543807920299116, 0.6930677925924869, 0.7184402241813643, 0.7048905554931792, 0.6537499962456508, 0.6173765319759497)
$colO = @(6.156959173689216, 6.21913647344897, 6.260920216710957, 6.278853745900051, 6.281886316330286, 6.26115840574019, 6.177648989203192, 6.042532138290426, 5.960764137482414, 5.927374971405129, 5.915279626158764, 5.917860170377651, 5.926368884056558, 5.931652162199299, 5.963022891038321, 5.983243645558304, 5.995232415869651, 5.999353123395451, 5.981054019728589, 5.923854777856889, 5.889668719607982, 5.907621645599193, 5.982042817046363, 6.042532138290426)

$cols = @('C', 'D', 'E', 'F', 'G', 'J', 'K', 'L', 'M', 'O')
$data = @{
    'C' = $colC
    'D' = $colD
    'E' = $colE
    'F' = $colF
    'G' = $colG
    'J' = $colJ
    'K' = $colK
    'L' = $colL
    'M' = $colM
    'O' = $colO
}

for ($r = 2; $r -le 25; $r++) {
    foreach ($c in $cols) {
        $idx = $r - 2
        $ws.Range("$c$r").Value = $data[$c][$idx]
    }
}
